$d = $word.ActiveDocument

# The last paragraph ("LDAP authentication.") ends with a trailing
# _GoBack bookmark (marks the last edit position). We need to insert a
# brand-new list paragraph after it, with the bookmark following along
# to the end of that new paragraph - exactly what happens in Word when
# you click at the very end of the document and press Enter, then type.
#
# Find/Replace alone keeps the trailing bookmark pinned to the true end
# of the replaced text (confirmed below), so first splice in the new
# text using a harmless marker in place of the paragraph break, then
# convert that marker into a real paragraph split via a Range op. The
# bookmark, already sitting after the new text, rides along with the
# paragraph-end as it's promoted to the new last paragraph.

$marker = "|||"
$newText = "Add functionality to add aliquots by scanning the barcode or CSV file into the layout."

$find1 = $d.Content
[void]$find1.Find.Execute("LDAP authentication.", $false, $false, $false, $false, $false, $true, 1, $false,
                     "LDAP authentication." + $marker + $newText, 2)

# Locate the marker and turn it into a genuine paragraph break, which
# splits the single paragraph into two (the second inheriting the same
# ListParagraph style / list numbering).
$find2 = $d.Content
[void]$find2.Find.Execute($marker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$find2.InsertParagraphAfter()

# Remove the marker text itself, leaving just the paragraph break.
$find2.Text = ""

$d.Save()
